$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "69.764.85"
$ws.Cells.Item(2, 5).Value = "  +0.66%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.429.45"
$ws.Cells.Item(3, 5).Value = "  +1.09%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "584.84"
$ws.Cells.Item(5, 5).Value = "  -0.65%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "177.00"
$ws.Cells.Item(6, 5).Value = "  -2.08%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "3.422.96"
$ws.Cells.Item(7, 5).Value = "  +1.03%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.01%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.59%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +2.10%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.585"
$ws.Cells.Item(11, 5).Value = "  -0.94%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "48.92"
$ws.Cells.Item(12, 5).Value = "  +0.31%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000284"
$ws.Cells.Item(13, 5).Value = "  -0.11%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "692.68"
$ws.Cells.Item(14, 5).Value = "  +1.67%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.975.88"
$ws.Cells.Item(15, 5).Value = "  +0.91%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "8.64"
$ws.Cells.Item(16, 5).Value = "  +0.06%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "69.774.50"
$ws.Cells.Item(17, 5).Value = "  +0.57%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "TRON"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.122"
$ws.Cells.Item(18, 5).Value = "  +1.08%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "WrappedEther"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "3.424.08"
$ws.Cells.Item(19, 5).Value = "  +0.76%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "17.68"
$ws.Cells.Item(20, 5).Value = "  -0.36%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.37"
$ws.Cells.Item(21, 5).Value = "  -0.23%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.898"
$ws.Cells.Item(22, 5).Value = "  -0.55%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.44"
$ws.Cells.Item(23, 5).Value = "  +0.27%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "16.94"
$ws.Cells.Item(24, 5).Value = "  -1.32%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "101.06"
$ws.Cells.Item(25, 5).Value = "  -2.47%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.92"
$ws.Cells.Item(26, 5).Value = "  -0.27%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -2.83%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "9.63"
$ws.Cells.Item(28, 5).Value = "  -0.29%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "33.50"
$ws.Cells.Item(29, 5).Value = "  -2.33%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "8.78"
$ws.Cells.Item(30, 5).Value = "  +0.39%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +1.81%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "573.01"
$ws.Cells.Item(32, 5).Value = "  +3.14%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.72"
$ws.Cells.Item(33, 5).Value = "  +0.10%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -1.92%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "OKB"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "58.36"
$ws.Cells.Item(35, 5).Value = "  +0.38%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.104"
$ws.Cells.Item(36, 5).Value = "  -2.62%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +0.15%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.570.16"
$ws.Cells.Item(38, 5).Value = "  -3.75%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -0.68%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "35.16"
$ws.Cells.Item(40, 5).Value = "  -1.61%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0₃0736"
$ws.Cells.Item(41, 5).Value = "  +3.91%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +0.37%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.69"
$ws.Cells.Item(43, 5).Value = "  +0.16%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +3.47%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "VeChain"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0420"
$ws.Cells.Item(45, 5).Value = "  -0.88%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "TheGraph"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.333"
$ws.Cells.Item(46, 5).Value = "  -2.04%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +3.81%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.65"
$ws.Cells.Item(48, 5).Value = "  -0.48%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.129"
$ws.Cells.Item(49, 5).Value = "  -1.09%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -0.23%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "133.02"
$ws.Cells.Item(51, 5).Value = "  +1.02%  "

Write-Output "Applied cryptos update"